$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 5.423951666666667
$ws.Range("H2").Value = 16.271855
$ws.Range("I2").Value = 0.4774188439413272
$ws.Range("J2").Value = 0.4774188439413271
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 218.7785543333333
$ws.Range("N2").Value = 656.3356630000001
$ws.Range("O2").Value = 0.7837094150017259
$ws.Range("P2").Value = 0.7837094150017259
$ws.Range("Q2").Value = 1186.644304407207
$ws.Range("R2").Value = 10679.79873966487
$ws.Range("S2").Value = 0.3741576428960578
$ws.Range("T2").Value = 0.3741576428960577

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 5.423951666666667
$ws.Range("H3").Value = 16.271855
$ws.Range("I3").Value = 0.4774188439413272
$ws.Range("J3").Value = 0.4774188439413271
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 46.29469433333333
$ws.Range("N3").Value = 138.884083
$ws.Range("O3").Value = 0.1658370397602197
$ws.Range("P3").Value = 0.1658370397602197
$ws.Range("Q3").Value = 251.1001844871072
$ws.Range("R3").Value = 2259.901660383965
$ws.Range("S3").Value = 0.07917372780497599
$ws.Range("T3").Value = 0.07917372780497596

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 5.423951666666667
$ws.Range("H4").Value = 16.271855
$ws.Range("I4").Value = 0.4774188439413272
$ws.Range("J4").Value = 0.4774188439413271
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 6.978882
$ws.Range("N4").Value = 20.936646
$ws.Range("O4").Value = 0.02499977909741928
$ws.Range("P4").Value = 0.02499977909741927
$ws.Range("Q4").Value = 37.85311865537
$ws.Range("R4").Value = 340.67806789833
$ws.Range("S4").Value = 0.01193536563547847
$ws.Range("T4").Value = 0.01193536563547846

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 5.423951666666667
$ws.Range("H5").Value = 16.271855
$ws.Range("I5").Value = 0.4774188439413272
$ws.Range("J5").Value = 0.4774188439413271
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 7.105616
$ws.Range("N5").Value = 21.316848
$ws.Range("O5").Value = 0.02545376614063513
$ws.Range("P5").Value = 0.02545376614063513
$ws.Range("Q5").Value = 38.54051774589333
$ws.Range("R5").Value = 346.86465971304
$ws.Range("S5").Value = 0.01215210760481492
$ws.Range("T5").Value = 0.01215210760481492

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.583504333333333
$ws.Range("H6").Value = 4.750513
$ws.Range("I6").Value = 0.1393808158066948
$ws.Range("J6").Value = 0.1393808158066948
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 218.7785543333333
$ws.Range("N6").Value = 656.3356630000001
$ws.Range("O6").Value = 0.7837094150017259
$ws.Range("P6").Value = 0.7837094150017259
$ws.Range("Q6").Value = 346.4367888272354
$ws.Range("R6").Value = 3117.931099445119
$ws.Range("S6").Value = 0.1092340576183281
$ws.Range("T6").Value = 0.1092340576183281

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.583504333333333
$ws.Range("H7").Value = 4.750513
$ws.Range("I7").Value = 0.1393808158066948
$ws.Range("J7").Value = 0.1393808158066948
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 46.29469433333333
$ws.Range("N7").Value = 138.884083
$ws.Range("O7").Value = 0.1658370397602197
$ws.Range("P7").Value = 0.1658370397602197
$ws.Range("Q7").Value = 73.30784908717544
$ws.Range("R7").Value = 659.770641784579
$ws.Range("S7").Value = 0.0231145018927467
$ws.Range("T7").Value = 0.0231145018927467

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1.583504333333333
$ws.Range("H8").Value = 4.750513
$ws.Range("I8").Value = 0.1393808158066948
$ws.Range("J8").Value = 0.1393808158066948
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 6.978882
$ws.Range("N8").Value = 20.936646
$ws.Range("O8").Value = 0.02499977909741928
$ws.Range("P8").Value = 0.02499977909741927
$ws.Range("Q8").Value = 11.051089888822
$ws.Range("R8").Value = 99.45980899939799
$ws.Range("S8").Value = 0.003484489605585455
$ws.Range("T8").Value = 0.003484489605585454

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1.583504333333333
$ws.Range("H9").Value = 4.750513
$ws.Range("I9").Value = 0.1393808158066948
$ws.Range("J9").Value = 0.1393808158066948
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 7.105616
$ws.Range("N9").Value = 21.316848
$ws.Range("O9").Value = 0.02545376614063513
$ws.Range("P9").Value = 0.02545376614063513
$ws.Range("Q9").Value = 11.25177372700267
$ws.Range("R9").Value = 101.265963543024
$ws.Range("S9").Value = 0.00354776669003455
$ws.Range("T9").Value = 0.003547766690034549

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.6660723333333333
$ws.Range("H10").Value = 1.998217
$ws.Range("I10").Value = 0.05862800830537802
$ws.Range("J10").Value = 0.05862800830537802
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 218.7785543333333
$ws.Range("N10").Value = 656.3356630000001
$ws.Range("O10").Value = 0.7837094150017259
$ws.Range("P10").Value = 0.7837094150017259
$ws.Range("Q10").Value = 145.7223421680968
$ws.Range("R10").Value = 1311.501079512871
$ws.Range("S10").Value = 0.04594732209172413
$ws.Range("T10").Value = 0.04594732209172413

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.6660723333333333
$ws.Range("H11").Value = 1.998217
$ws.Range("I11").Value = 0.05862800830537802
$ws.Range("J11").Value = 0.05862800830537802
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 46.29469433333333
$ws.Range("N11").Value = 138.884083
$ws.Range("O11").Value = 0.1658370397602197
$ws.Range("P11").Value = 0.1658370397602197
$ws.Range("Q11").Value = 30.83561507555677
$ws.Range("R11").Value = 277.520535680011
$ws.Range("S11").Value = 0.009722695344401464
$ws.Range("T11").Value = 0.009722695344401464

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.6660723333333333
$ws.Range("H12").Value = 1.998217
$ws.Range("I12").Value = 0.05862800830537802
$ws.Range("J12").Value = 0.05862800830537802
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 6.978882
$ws.Range("N12").Value = 20.936646
$ws.Range("O12").Value = 0.02499977909741928
$ws.Range("P12").Value = 0.02499977909741927
$ws.Range("Q12").Value = 4.648440217797999
$ws.Range("R12").Value = 41.835961960182
$ws.Range("S12").Value = 0.001465687256556113
$ws.Range("T12").Value = 0.001465687256556113

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.6660723333333333
$ws.Range("H13").Value = 1.998217
$ws.Range("I13").Value = 0.05862800830537802
$ws.Range("J13").Value = 0.05862800830537802
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 7.105616
$ws.Range("N13").Value = 21.316848
$ws.Range("O13").Value = 0.02545376614063513
$ws.Range("P13").Value = 0.02545376614063513
$ws.Range("Q13").Value = 4.732854228890666
$ws.Range("R13").Value = 42.595688060016
$ws.Range("S13").Value = 0.001492303612696306
$ws.Range("T13").Value = 0.001492303612696306

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 3.687463666666666
$ws.Range("H14").Value = 11.062391
$ws.Range("I14").Value = 0.3245723319466
$ws.Range("J14").Value = 0.3245723319466
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 218.7785543333333
$ws.Range("N14").Value = 656.3356630000001
$ws.Range("O14").Value = 0.7837094150017259
$ws.Range("P14").Value = 0.7837094150017259
$ws.Range("Q14").Value = 806.737970150026
$ws.Range("R14").Value = 7260.641731350233
$ws.Range("S14").Value = 0.2543703923956158
$ws.Range("T14").Value = 0.2543703923956158

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 3.687463666666666
$ws.Range("H15").Value = 11.062391
$ws.Range("I15").Value = 0.3245723319466
$ws.Range("J15").Value = 0.3245723319466
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 46.29469433333333
$ws.Range("N15").Value = 138.884083
$ws.Range("O15").Value = 0.1658370397602197
$ws.Range("P15").Value = 0.1658370397602197
$ws.Range("Q15").Value = 170.7100033136059
$ws.Range("R15").Value = 1536.390029822453
$ws.Range("S15").Value = 0.05382611471809551
$ws.Range("T15").Value = 0.05382611471809551

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 3.687463666666666
$ws.Range("H16").Value = 11.062391
$ws.Range("I16").Value = 0.3245723319466
$ws.Range("J16").Value = 0.3245723319466
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 6.978882
$ws.Range("N16").Value = 20.936646
$ws.Range("O16").Value = 0.02499977909741928
$ws.Range("P16").Value = 0.02499977909741927
$ws.Range("Q16").Value = 25.734373808954
$ws.Range("R16").Value = 231.609364280586
$ws.Range("S16").Value = 0.00811423659979924
$ws.Range("T16").Value = 0.00811423659979924

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 3.687463666666666
$ws.Range("H17").Value = 11.062391
$ws.Range("I17").Value = 0.3245723319466
$ws.Range("J17").Value = 0.3245723319466
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 7.105616
$ws.Range("N17").Value = 21.316848
$ws.Range("O17").Value = 0.02545376614063513
$ws.Range("P17").Value = 0.02545376614063513
$ws.Range("Q17").Value = 26.20170082928533
$ws.Range("R17").Value = 235.815307463568
$ws.Range("S17").Value = 0.008261588233089352
$ws.Range("T17").Value = 0.008261588233089352
